# Auto-generated edit script applying the Titan_Profits diff
# Values below were derived from the unified XML diff, mapped to their
# correct worksheet (ALC/ARM/BSM/CRP/CUL/LTW/WVR) by matching the unique
# 'Leve Item ID' (column G) context value recorded alongside each hunk.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 605.9375
$ws.Range("I19").Value = 493
$ws.Range("J19").Value = 673.7
$ws.Range("K19").Value = 493
$ws.Range("L19").Value = 673.7
$ws.Range("M19").Value = -318
$ws.Range("N19").Value = -1023.7
$ws.Range("H106").Value = 27780056
$ws.Range("I106").Value = 27780056
$ws.Range("J106").Value = 0
$ws.Range("K106").Value = 27780056
$ws.Range("L106").Value = 0
$ws.Range("M106").Value = ""
$ws.Range("N106").Value = -27779425
$ws.Range("H111").Value = 1610.8182
$ws.Range("I111").Value = 1275.5714
$ws.Range("J111").Value = 2197.5
$ws.Range("K111").Value = 3826.7142
$ws.Range("L111").Value = 6592.5
$ws.Range("M111").Value = -759.7142000000003
$ws.Range("N111").Value = -12726.5
$ws.Range("H129").Value = 1201.5294
$ws.Range("J129").Value = 1434.6154
$ws.Range("L129").Value = 4303.8462
$ws.Range("N129").Value = -14303.8462
$ws.Range("H132").Value = 24364.71
$ws.Range("I132").Value = 24916.045
$ws.Range("J132").Value = 106
$ws.Range("K132").Value = 74748.13499999999
$ws.Range("L132").Value = 318
$ws.Range("M132").Value = -72218.13499999999
$ws.Range("N132").Value = -5378
$ws.Range("H135").Value = 1630.6923
$ws.Range("I135").Value = 1969.3
$ws.Range("J135").Value = 502
$ws.Range("K135").Value = 17723.7
$ws.Range("L135").Value = 4518
$ws.Range("M135").Value = -15188.7
$ws.Range("N135").Value = -9588

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 5527.0347
$ws.Range("I74").Value = 1337.1
$ws.Range("J74").Value = 14838
$ws.Range("K74").Value = 1337.1
$ws.Range("L74").Value = 14838
$ws.Range("M74").Value = -463.0999999999999
$ws.Range("N74").Value = -16586
$ws.Range("H77").Value = 5527.0347
$ws.Range("I77").Value = 1337.1
$ws.Range("J77").Value = 14838
$ws.Range("K77").Value = 6685.5
$ws.Range("L77").Value = 74190
$ws.Range("M77").Value = -2317.5
$ws.Range("N77").Value = -82926
$ws.Range("H132").Value = 2153.0527
$ws.Range("I132").Value = 1824.3846
$ws.Range("K132").Value = 5473.1538
$ws.Range("M132").Value = -2943.1538

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1508.2307
$ws.Range("I86").Value = 1567.8334
$ws.Range("J86").Value = 1457.1428
$ws.Range("K86").Value = 1567.8334
$ws.Range("L86").Value = 1457.1428
$ws.Range("M86").Value = -444.8334
$ws.Range("N86").Value = -3703.1428
$ws.Range("H89").Value = 1508.2307
$ws.Range("I89").Value = 1567.8334
$ws.Range("J89").Value = 1457.1428
$ws.Range("K89").Value = 7839.166999999999
$ws.Range("L89").Value = 7285.714
$ws.Range("M89").Value = -2223.166999999999
$ws.Range("N89").Value = -18517.714
$ws.Range("H134").Value = 30306548
$ws.Range("I134").Value = 45457156
$ws.Range("J134").Value = 5332.909
$ws.Range("K134").Value = 136371468
$ws.Range("L134").Value = 15998.727
$ws.Range("M134").Value = -136368933
$ws.Range("N134").Value = -21068.727

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H11").Value = 900
$ws.Range("I11").Value = 900
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 900
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = ""
$ws.Range("N11").Value = -760
$ws.Range("H41").Value = 14296
$ws.Range("J41").Value = 38888
$ws.Range("L41").Value = 38888
$ws.Range("N41").Value = -39744
$ws.Range("H50").Value = 17500
$ws.Range("I50").Value = 8333.333000000001
$ws.Range("J50").Value = 26666.666
$ws.Range("K50").Value = 8333.333000000001
$ws.Range("L50").Value = 26666.666
$ws.Range("M50").Value = -7708.333000000001
$ws.Range("N50").Value = -27916.666
$ws.Range("H51").Value = 95237.57000000001
$ws.Range("J51").Value = 95237.57000000001
$ws.Range("L51").Value = 95237.57000000001
$ws.Range("N51").Value = -96709.57000000001
$ws.Range("J59").Value = 20000
$ws.Range("L59").Value = 20000
$ws.Range("N59").Value = -22290
$ws.Range("H61").Value = 95237.57000000001
$ws.Range("J61").Value = 95237.57000000001
$ws.Range("L61").Value = 95237.57000000001
$ws.Range("N61").Value = -95933.57000000001
$ws.Range("H62").Value = 19094.6
$ws.Range("I62").Value = 23719.908
$ws.Range("J62").Value = 6375
$ws.Range("K62").Value = 23719.908
$ws.Range("L62").Value = 6375
$ws.Range("M62").Value = -23095.908
$ws.Range("N62").Value = -7623
$ws.Range("H65").Value = 19094.6
$ws.Range("I65").Value = 23719.908
$ws.Range("J65").Value = 6375
$ws.Range("K65").Value = 118599.54
$ws.Range("L65").Value = 31875
$ws.Range("M65").Value = -115479.54
$ws.Range("N65").Value = -38115
$ws.Range("H68").Value = 24095
$ws.Range("I68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("M68").Value = ""
$ws.Range("H71").Value = 24095
$ws.Range("I71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("M71").Value = ""
$ws.Range("H74").Value = 18576.4
$ws.Range("J74").Value = 18576.4
$ws.Range("L74").Value = 18576.4
$ws.Range("N74").Value = -20324.4
$ws.Range("H77").Value = 18576.4
$ws.Range("J77").Value = 18576.4
$ws.Range("L77").Value = 55729.2
$ws.Range("N77").Value = -64465.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 6250
$ws.Range("I56").Value = 6250
$ws.Range("K56").Value = 6250
$ws.Range("M56").Value = -5720
$ws.Range("H131").Value = 12347746
$ws.Range("I131").Value = 475
$ws.Range("J131").Value = 14495098
$ws.Range("K131").Value = 1425
$ws.Range("L131").Value = 43485294
$ws.Range("M131").Value = 3615
$ws.Range("N131").Value = -43495374

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1224.3334
$ws.Range("J46").Value = 1559.8
$ws.Range("L46").Value = 1559.8
$ws.Range("N46").Value = -1935.8
$ws.Range("H104").Value = 30000
$ws.Range("J104").Value = 30000
$ws.Range("L104").Value = 30000
$ws.Range("N104").Value = -36988

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H104").Value = 29900
$ws.Range("J104").Value = 29900
$ws.Range("L104").Value = 29900
$ws.Range("N104").Value = -36888
$ws.Range("H122").Value = 2115.8823
$ws.Range("I122").Value = 1984.1111
$ws.Range("J122").Value = 2264.125
$ws.Range("K122").Value = 5952.3333
$ws.Range("L122").Value = 6792.375
$ws.Range("M122").Value = -3502.3333
$ws.Range("N122").Value = -11692.375
$ws.Range("H132").Value = 1750.6735
$ws.Range("I132").Value = 1479.2778
$ws.Range("J132").Value = 2502.2307
$ws.Range("K132").Value = 4437.8334
$ws.Range("L132").Value = 7506.6921
$ws.Range("M132").Value = -1907.8334
$ws.Range("N132").Value = -12566.6921
$ws.Range("H137").Value = 61119.168
$ws.Range("J137").Value = 97238.336
$ws.Range("L137").Value = 97238.336
$ws.Range("N137").Value = -107438.336
$ws.Range("H140").Value = 72467.73
$ws.Range("J140").Value = 72467.73
$ws.Range("L140").Value = 72467.73
$ws.Range("N140").Value = -82827.73

